$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "IPPIS No: {ippisNumber}" -> "IPPIS No: {ippis}"  (first table)
# ------------------------------------------------------------------
$tbl1 = $d.Tables.Item(1)
$tbl1.Range.Find.Execute("ippisNumber", $true, $false, $false, $false, $false, `
    $true, 1, $false, "ippis", 2)

# ------------------------------------------------------------------
# 2) "{ippisNumber}" -> "{ippis" + _GoBack bookmark + "}"  (second table)
# ------------------------------------------------------------------
$tbl2 = $d.Tables.Item(2)
$hit = $tbl2.Range.Duplicate
$hit.Find.Execute("{ippisNumber}")
$hitStart = $hit.Start
$hitEnd = $hit.End

$replaceRange = $d.Range($hitStart, $hitEnd)
$replaceRange.Text = "{ippis}"

$splitPos = $hitStart + 6
$bookmarkRange = $d.Range($splitPos, $splitPos)

# ------------------------------------------------------------------
# 3) Move the _GoBack bookmark from the end of the document to the
#    split point created above.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# ------------------------------------------------------------------
# 4) Resize the second table's grid columns.
# ------------------------------------------------------------------
$tbl2.Columns.Item(1).Width = 151.85
$tbl2.Columns.Item(2).Width = 56.45
$tbl2.Columns.Item(3).Width = 104.7
$tbl2.Columns.Item(4).Width = 119.8
